# Actualización automática hashcode mar may  7 01:40:46 CEST 2019
# Updates the hashcode values (column B) for a set of rows in the
# "hashcode.csv" sheet to reflect newly computed MD5 hashes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    11  = "ce3cbe9e64f802728b07a273f55120a3"
    15  = "ca97fbdc43e8366a17b0d9dad8e5c063"
    29  = "ca4e7eb493cd0a8237bdd4e609b357c8"
    121 = "8012a2d944ac783cf477cc4c54fdc599"
    169 = "0d8770a9fc02564072eebdd04c43a5e5"
    213 = "94bd705449ea5f0d10f94f453af6e990"
    281 = "624f64092c79ea3bd74ebec9e1b940b3"
    339 = "bc143bd69e4e945d769f1abedefd36f1"
    461 = "1ef706b6c3ed1200c4cf9b3b06993c8f"
    500 = "0d4ef51694dce6eec8c2c56e296b0feb"
    502 = "7a1bcfc6aebbf91920e3541b8fdcbe95"
    517 = "bdd041522b328e847d0665f3fda436b3"
    547 = "96159b90fc80a73cc6b204aedf87156c"
    616 = "ca7a419aae4ba29d65207207f9eca58d"
    627 = "e05231403251ca69c2359ff132eb8959"
    629 = "69b3b7d4dd76f850536665fd29743f8f"
    655 = "82747305e49304ffb6401a0786ba856f"
    665 = "946edc2d5916c25b11b9d997aea52506"
    685 = "b72723fe76e241d54eaf584cfe9f49e0"
    733 = "f7b1ee1975b192822078cf9e38a17f72"
    819 = "2235951dc2d4e550c50e4e31bba45850"
    862 = "937e5eec33ddae54a59b9fe523ed70bf"
    874 = "c9c849f03081bb7a17b5eba5feebb7ea"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
